# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '89.528.89'
$ws.Range('E2').Value = '  +4.80%  '
$ws.Range('D3').Value = "'" + '3.232.34'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('D4').Value = "'" + '1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'" + '214.26'
$ws.Range('E5').Value = '  +4.06%  '
$ws.Range('D6').Value = "'" + '621.02'
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('D7').Value = "'" + '0.404'
$ws.Range('E7').Value = '  +16.83%  '
$ws.Range('D8').Value = "'" + '0.706'
$ws.Range('E8').Value = '  +11.49%  '
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').Value = "'" + '3.224.72'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('D11').Value = "'" + '0.562'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = "'" + '0.179'
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('E13').Value = '  +4.93%  '
$ws.Range('D14').Value = "'" + '5.38'
$ws.Range('E14').Value = '  +4.01%  '
$ws.Range('D15').Value = "'" + '3.810.86'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = "'" + '33.36'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = "'" + '89.116.95'
$ws.Range('E17').Value = '  +4.39%  '
$ws.Range('D18').Value = "'" + '3.218.08'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  +11.44%  '
$ws.Range('D20').Value = "'" + '13.86'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = "'" + '423.35'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').Value = "'" + '8.70'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('E23').Value = '  +1.38%  '
$ws.Range('D24').Value = "'" + '0.0000164'
$ws.Range('E24').Value = '  +31.90%  '
$ws.Range('D25').Value = "'" + '5.43'
$ws.Range('E25').Value = '  +7.58%  '
$ws.Range('D26').Value = "'" + '12.57'
$ws.Range('E26').Value = '  +3.19%  '
$ws.Range('D27').Value = "'" + '3.333.64'
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('D28').Value = "'" + '74.98'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D30').Value = "'" + '0.175'
$ws.Range('E30').Value = '  +4.00%  '
$ws.Range('D31').Value = "'" + '0.998'
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').Value = "'" + '565.33'
$ws.Range('E32').Value = '  +5.75%  '
$ws.Range('D33').Value = "'" + '8.44'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = "'" + '7.14'
$ws.Range('E34').Value = '  +8.82%  '
$ws.Range('E35').Value = '  -4.10%  '
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('D37').Value = "'" + '0.134'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('E38').Value = '  +1.58%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = "'" + '21.89'
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = "'" + '3.21'
$ws.Range('E40').Value = '  +13.38%  '
$ws.Range('D41').Value = "'" + '1.00'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').Value = "'" + '1.97'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('D43').Value = "'" + '0.384'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = "'" + '151.83'
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('D46').Value = "'" + '177.98'
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').Value = "'" + '43.65'
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('E48').Value = '  +10.12%  '
$ws.Range('D49').Value = "'" + '1.27'
$ws.Range('E49').Value = '  -1.54%  '
$ws.Range('D50').Value = "'" + '24.72'
$ws.Range('E50').Value = '  +4.86%  '
$ws.Range('D51').Value = "'" + '4.05'
$ws.Range('E51').Value = '  -2.71%  '
